$d = $word.ActiveDocument
$bullet = [char]0x2022

# ---------------------------------------------------------------------------
# 0. Remove the stray "_GoBack" bookmark that currently sits at the very end
#    of the document, near "Fluent in Spanish". It is being relocated to the
#    new OBJECTIVE paragraph we add further down, so get rid of the old one
#    first (before any new bookmark with the same name gets created).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 1. Contact-info line: merge the (previously spell-check-split) runs into a
#    single run and append the github / linkedin handles.
# ---------------------------------------------------------------------------
$contactPara = $d.Paragraphs.Item(3)
$contactRange = $contactPara.Range
$contactRange.End = $contactRange.End - 1
$contactRange.Text = "PLACEHOLDER_CONTACT"
$contactRange2 = $d.Paragraphs.Item(3).Range
$contactRange2.End = $contactRange2.End - 1
$contactRange2.Text = "(803)389-6750 $bullet danielmartincraig@gmail.com $bullet github.com/danielmartincraig $bullet linkedin.com/danielcraig23"

# ---------------------------------------------------------------------------
# 2. Insert a new "OBJECTIVE:" paragraph right after the contact-info line
#    (i.e. right before the "EDUCATION:" heading, which is paragraph 4).
# ---------------------------------------------------------------------------
$eduPara = $d.Paragraphs.Item(4)
$eduRange = $eduPara.Range
$eduRange.InsertParagraphBefore()

$objPara = $d.Paragraphs.Item(4)
$objRange = $objPara.Range

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$objXml = "<w:p $wns>" +
  "<w:pPr><w:pStyle w:val=`"Heading1`"/><w:rPr><w:sz w:val=`"24`"/></w:rPr></w:pPr>" +
  "<w:r><w:t xml:space=`"preserve`">OBJECTIVE: </w:t></w:r>" +
  "<w:r><w:rPr><w:sz w:val=`"24`"/></w:rPr><w:t xml:space=`"preserve`">Eager to drive back-end solutions at </w:t></w:r>" +
  "<w:r><w:rPr><w:sz w:val=`"24`"/></w:rPr><w:t>Owens Corning</w:t></w:r>" +
  "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/>" +
  "<w:r><w:rPr><w:sz w:val=`"24`"/></w:rPr><w:t xml:space=`"preserve`"> on a full-time basis</w:t></w:r>" +
  "</w:p>"

$objRange.InsertXML($objXml)

# ---------------------------------------------------------------------------
# 3. "Web Engineering I and II" bullet: merge the grammar-check-split runs
#    into a single run.
# ---------------------------------------------------------------------------
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "$bullet    Web Engineering I and II`r") {
        $targetPara = $para
        break
    }
}
if ($targetPara -ne $null) {
    $webRange = $targetPara.Range
    $webRange.End = $webRange.End - 1
    $webRange.Text = "PLACEHOLDER_WEB"
    $webPara2 = $null
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $para = $d.Paragraphs.Item($i)
        if ($para.Range.Text -eq "PLACEHOLDER_WEB`r") {
            $webPara2 = $para
            break
        }
    }
    $webRange2 = $webPara2.Range
    $webRange2.End = $webRange2.End - 1
    $webRange2.Text = "$bullet    Web Engineering I and II"
}

Write-Host "done"
